$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.353.75"
$ws.Range("E2").Value = "  +3.70%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.245.94"
$ws.Range("E3").Value = "  +2.53%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.16"
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.52"
$ws.Range("E6").Value = "  +6.01%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -4.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.244.60"
$ws.Range("E9").Value = "  +2.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.132"
$ws.Range("E10").Value = "  +5.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.77"
$ws.Range("E11").Value = "  +3.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.414"
$ws.Range("E12").Value = "  +4.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.810.18"
$ws.Range("E13").Value = "  +2.53%  "
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.49"
$ws.Range("E15").Value = "  +4.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.331.78"
$ws.Range("E16").Value = "  +3.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000168"
$ws.Range("E17").Value = "  +2.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.251.83"
$ws.Range("E18").Value = "  +2.84%  "
$ws.Range("E19").Value = "  +1.75%  "
$ws.Range("E20").Value = "  +5.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.72"
$ws.Range("E21").Value = "  +5.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.62"
$ws.Range("E22").Value = "  +4.79%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.31"
$ws.Range("E24").Value = "  +3.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.511"
$ws.Range("E25").Value = "  +2.17%  "
$ws.Range("E26").Value = "  +1.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.59"
$ws.Range("E27").Value = "  -1.12%  "
$ws.Range("E28").Value = "  +2.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.80"
$ws.Range("E30").Value = "  +7.80%  "
$ws.Range("E31").Value = "  +3.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.67"
$ws.Range("E32").Value = "  +2.97%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  +5.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.93"
$ws.Range("E35").Value = "  +3.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.99"
$ws.Range("E36").Value = "  +6.19%  "
$ws.Range("E37").Value = "  +3.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.848"
$ws.Range("E38").Value = "  +1.27%  "
$ws.Range("E39").Value = "  +4.78%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.83"
$ws.Range("E40").Value = "  +13.36%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.86"
$ws.Range("E41").Value = "  +2.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.60"
$ws.Range("E42").Value = "  +9.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.60"
$ws.Range("E43").Value = "  +3.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "356.73"
$ws.Range("E44").Value = "  +8.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.726.14"
$ws.Range("E45").Value = "  +2.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.36"
$ws.Range("E46").Value = "  +4.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.80"
$ws.Range("E47").Value = "  +3.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0678"
$ws.Range("E48").Value = "  +2.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0279"
$ws.Range("E49").Value = "  +1.37%  "
$ws.Range("E50").Value = "  +5.83%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.102"
$ws.Range("E51").Value = "  -0.99%  "
